# Updated cryptos list on Mon May  8 07:43:45 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns of the crypto
# table with the latest scraped figures, and corrects the order of the
# RenderToken / PaxDollar rows (49-50).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores numeric-looking labels as literal text (e.g. "28.217.52",
# "1.004", "0.000009019") rather than real numbers, so force a Text number
# format on every D cell we are about to rewrite. Without this, Excel would
# silently reinterpret strings such as "1.005" or "0.000009008" as floating
# point numbers and corrupt their original formatting/precision.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Now write the refreshed values. Column D is assigned via .Formula (the
# NumberFormat="@" set above makes the engine keep it as literal text);
# columns B, C and E are plain text already, so .Value is sufficient.
$ws.Range("D2").Formula = '28.201.34'
$ws.Range("E2").Value = '  -2.71%  '
$ws.Range("D3").Formula = '1.864.91'
$ws.Range("E3").Value = '  -2.38%  '
$ws.Range("D4").Formula = '1.005'
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Formula = '318.57'
$ws.Range("E5").Value = '  -1.91%  '
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("D7").Formula = '0.4392'
$ws.Range("E7").Value = '  -4.33%  '
$ws.Range("D8").Formula = '0.3683'
$ws.Range("E8").Value = '  -3.74%  '
$ws.Range("D9").Formula = '0.07484'
$ws.Range("E9").Value = '  -3.02%  '
$ws.Range("D10").Formula = '0.9330'
$ws.Range("E10").Value = '  -4.85%  '
$ws.Range("D11").Formula = '21.26'
$ws.Range("E11").Value = '  -3.78%  '
$ws.Range("D12").Formula = '1.862.11'
$ws.Range("E12").Value = '  -2.36%  '
$ws.Range("D13").Formula = '6.684'
$ws.Range("E13").Value = '  -3.64%  '
$ws.Range("D14").Formula = '5.425'
$ws.Range("E14").Value = '  -4.43%  '
$ws.Range("D15").Formula = '0.06899'
$ws.Range("E15").Value = '  -1.71%  '
$ws.Range("D16").Formula = '1.005'
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").Formula = '81.66'
$ws.Range("E17").Value = '  -2.79%  '
$ws.Range("D18").Formula = '0.000009008'
$ws.Range("E18").Value = '  -4.94%  '
$ws.Range("D19").Formula = '1.004'
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("D20").Formula = '15.85'
$ws.Range("E20").Value = '  -5.26%  '
$ws.Range("D21").Formula = '28.176.40'
$ws.Range("E21").Value = '  -2.73%  '
$ws.Range("D22").Formula = '5.105'
$ws.Range("E22").Value = '  -4.06%  '
$ws.Range("E23").Value = '  -1.04%  '
$ws.Range("D24").Formula = '2.155.95'
$ws.Range("E24").Value = '  +1.30%  '
$ws.Range("D25").Formula = '2.012'
$ws.Range("E25").Value = '  -3.77%  '
$ws.Range("D26").Formula = '155.10'
$ws.Range("E26").Value = '  -2.07%  '
$ws.Range("D27").Formula = '18.36'
$ws.Range("E27").Value = '  -3.59%  '
$ws.Range("D28").Formula = '5.311'
$ws.Range("E28").Value = '  -6.35%  '
$ws.Range("D29").Formula = '113.02'
$ws.Range("E29").Value = '  -3.85%  '
$ws.Range("D30").Formula = '1.718'
$ws.Range("E30").Value = '  -7.18%  '
$ws.Range("E31").Value = '  -2.89%  '
$ws.Range("D32").Formula = '4.830'
$ws.Range("E32").Value = '  -4.80%  '
$ws.Range("D33").Formula = '0.7904'
$ws.Range("E33").Value = '  -8.84%  '
$ws.Range("D34").Formula = '1.171'
$ws.Range("E34").Value = '  -6.45%  '
$ws.Range("D35").Formula = '2.942'
$ws.Range("E35").Value = '  -2.72%  '
$ws.Range("E36").Value = '  +0.18%  '
$ws.Range("D37").Formula = '1.123'
$ws.Range("E37").Value = '  -2.86%  '
$ws.Range("D38").Formula = '0.05426'
$ws.Range("E38").Value = '  -5.60%  '
$ws.Range("D39").Formula = '0.01965'
$ws.Range("E39").Value = '  -3.79%  '
$ws.Range("D40").Formula = '2.951'
$ws.Range("E40").Value = '  +2.50%  '
$ws.Range("D41").Formula = '0.5244'
$ws.Range("E41").Value = '  -4.96%  '
$ws.Range("D42").Formula = '6.972'
$ws.Range("E42").Value = '  -6.18%  '
$ws.Range("D43").Formula = '0.1674'
$ws.Range("E43").Value = '  -4.75%  '
$ws.Range("D44").Formula = '8.681'
$ws.Range("E44").Value = '  -6.92%  '
$ws.Range("D45").Formula = '0.06737'
$ws.Range("E45").Value = '  -1.67%  '
$ws.Range("D46").Formula = '0.4860'
$ws.Range("E46").Value = '  -6.33%  '
$ws.Range("D47").Formula = '10.51'
$ws.Range("E47").Value = '  -7.06%  '
$ws.Range("D48").Formula = '106.92'
$ws.Range("E48").Value = '  -3.73%  '
$ws.Range("B49").Value = 'PaxDollar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D49").Formula = '1.003'
$ws.Range("E49").Value = '  +0.17%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Formula = '1.909'
$ws.Range("E50").Value = '  -7.54%  '
$ws.Range("D51").Formula = '1.668'
$ws.Range("E51").Value = '  -6.45%  '
